$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 449.0909
$ws.Range("I6").Value = 535.6667
$ws.Range("J6").Value = 59.5
$ws.Range("K6").Value = 1607.0001
$ws.Range("L6").Value = 178.5
$ws.Range("M6").Value = -1495.0001
$ws.Range("N6").Value = -402.5
$ws.Range("H15").Value = 1166.4722
$ws.Range("I15").Value = 1166.4722
$ws.Range("K15").Value = 3499.4166
$ws.Range("M15").Value = -3330.4166
$ws.Range("H32").Value = 11683.167
$ws.Range("I32").Value = 1999.6666
$ws.Range("J32").Value = 21366.666
$ws.Range("K32").Value = 1999.6666
$ws.Range("L32").Value = 21366.666
$ws.Range("M32").Value = -1673.6666
$ws.Range("N32").Value = -22018.666
$ws.Range("H70").Value = 168164
$ws.Range("J70").Value = 501249.5
$ws.Range("L70").Value = 1503748.5
$ws.Range("N70").Value = -1504288.5
$ws.Range("H73").Value = 168164
$ws.Range("J73").Value = 501249.5
$ws.Range("L73").Value = 1503748.5
$ws.Range("N73").Value = -1505620.5
$ws.Range("H98").Value = 3445.9666
$ws.Range("I98").Value = 1039.65
$ws.Range("K98").Value = 1039.65
$ws.Range("M98").Value = 458.3499999999999
$ws.Range("H100").Value = 4357
$ws.Range("I100").Value = 3844
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 3844
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -3303
$ws.Range("N100").Value = -11082
$ws.Range("H111").Value = 999.5
$ws.Range("I111").Value = 999.5
$ws.Range("K111").Value = 2998.5
$ws.Range("M111").Value = 68.5
$ws.Range("H122").Value = 3445.9666
$ws.Range("I122").Value = 1039.65
$ws.Range("K122").Value = 3118.95
$ws.Range("M122").Value = -668.9500000000003
$ws.Range("H127").Value = 479.57144
$ws.Range("I127").Value = 479.57144
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1438.71432
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 3521.28568
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 1726.7894
$ws.Range("I129").Value = 1211.7858
$ws.Range("J129").Value = 3168.8
$ws.Range("K129").Value = 3635.3574
$ws.Range("L129").Value = 9506.400000000001
$ws.Range("M129").Value = 1364.6426
$ws.Range("N129").Value = -19506.4
$ws.Range("H132").Value = 1767.4286
$ws.Range("I132").Value = 1797.7764
$ws.Range("J132").Value = 1479.125
$ws.Range("K132").Value = 5393.3292
$ws.Range("L132").Value = 4437.375
$ws.Range("M132").Value = -2863.3292
$ws.Range("N132").Value = -9497.375
$ws.Range("H135").Value = 1273.6
$ws.Range("I135").Value = 1193.409
$ws.Range("J135").Value = 1861.6666
$ws.Range("K135").Value = 10740.681
$ws.Range("L135").Value = 16754.9994
$ws.Range("M135").Value = -8205.681
$ws.Range("N135").Value = -21824.9994
$ws.Range("H137").Value = 2563.875
$ws.Range("I137").Value = 2362.2
$ws.Range("J137").Value = 2900
$ws.Range("K137").Value = 7086.599999999999
$ws.Range("L137").Value = 8700
$ws.Range("M137").Value = -4536.599999999999
$ws.Range("N137").Value = -13800
$ws.Range("H138").Value = 2848.16
$ws.Range("I138").Value = 1179.7675
$ws.Range("J138").Value = 5090.0625
$ws.Range("K138").Value = 3539.3025
$ws.Range("L138").Value = 15270.1875
$ws.Range("M138").Value = 1600.6975
$ws.Range("N138").Value = -25550.1875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7106.1514
$ws.Range("I32").Value = 948.3333
$ws.Range("K32").Value = 948.3333
$ws.Range("M32").Value = -661.3333
$ws.Range("H43").Value = 15933.875
$ws.Range("J43").Value = 12096
$ws.Range("L43").Value = 12096
$ws.Range("N43").Value = -12722
$ws.Range("H45").Value = 2531.7058
$ws.Range("I45").Value = 2345.9666
$ws.Range("K45").Value = 2345.9666
$ws.Range("M45").Value = -1968.9666
$ws.Range("H53").Value = 25647.834
$ws.Range("I53").Value = 20777.6
$ws.Range("K53").Value = 20777.6
$ws.Range("M53").Value = -20095.6
$ws.Range("H61").Value = 7187.7646
$ws.Range("I61").Value = 4442.7144
$ws.Range("K61").Value = 4442.7144
$ws.Range("M61").Value = -4230.7144
$ws.Range("H74").Value = 1896.9656
$ws.Range("I74").Value = 1600.75
$ws.Range("K74").Value = 1600.75
$ws.Range("M74").Value = -726.75
$ws.Range("H77").Value = 1896.9656
$ws.Range("I77").Value = 1600.75
$ws.Range("K77").Value = 8003.75
$ws.Range("M77").Value = -3635.75
$ws.Range("H132").Value = 4388.5557
$ws.Range("I132").Value = 3856.7856
$ws.Range("J132").Value = 6249.75
$ws.Range("K132").Value = 11570.3568
$ws.Range("L132").Value = 18749.25
$ws.Range("M132").Value = -9040.356800000001
$ws.Range("N132").Value = -23809.25
$ws.Range("H136").Value = 7187.7646
$ws.Range("I136").Value = 4442.7144
$ws.Range("K136").Value = 13328.1432
$ws.Range("M136").Value = -10778.1432
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6242.706
$ws.Range("I31").Value = 5951.8945
$ws.Range("J31").Value = 6611.067
$ws.Range("K31").Value = 5951.8945
$ws.Range("L31").Value = 6611.067
$ws.Range("M31").Value = -5656.8945
$ws.Range("N31").Value = -7201.067
$ws.Range("H34").Value = 6242.706
$ws.Range("I34").Value = 5951.8945
$ws.Range("J34").Value = 6611.067
$ws.Range("K34").Value = 5951.8945
$ws.Range("L34").Value = 6611.067
$ws.Range("M34").Value = -5749.8945
$ws.Range("N34").Value = -7015.067
$ws.Range("H55").Value = 7500
$ws.Range("I55").Value = 7500
$ws.Range("K55").Value = 7500
$ws.Range("M55").Value = -7185
$ws.Range("H74").Value = 49999
$ws.Range("J74").Value = 49999
$ws.Range("L74").Value = 49999
$ws.Range("N74").Value = -51747
$ws.Range("H77").Value = 49999
$ws.Range("J77").Value = 49999
$ws.Range("L77").Value = 149997
$ws.Range("N77").Value = -158733
$ws.Range("H99").Value = 8394.125
$ws.Range("I99").Value = 9020.6
$ws.Range("K99").Value = 9020.6
$ws.Range("M99").Value = -7522.6
$ws.Range("H102").Value = 70241
$ws.Range("J102").Value = 70241
$ws.Range("L102").Value = 70241
$ws.Range("N102").Value = -75109
$ws.Range("H103").Value = 47999
$ws.Range("I103").Value = 36999
$ws.Range("K103").Value = 36999
$ws.Range("M103").Value = -35827
$ws.Range("H126").Value = 8394.125
$ws.Range("I126").Value = 9020.6
$ws.Range("K126").Value = 27061.8
$ws.Range("M126").Value = -24591.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 862.9167
$ws.Range("I26").Value = 1029.25
$ws.Range("J26").Value = 530.25
$ws.Range("K26").Value = 3087.75
$ws.Range("L26").Value = 1590.75
$ws.Range("M26").Value = -2799.75
$ws.Range("N26").Value = -2166.75
$ws.Range("H38").Value = 80.40000000000001
$ws.Range("I38").Value = 38.5
$ws.Range("J38").Value = 108.333336
$ws.Range("K38").Value = 115.5
$ws.Range("L38").Value = 325.000008
$ws.Range("M38").Value = 231.5
$ws.Range("N38").Value = -1019.000008
$ws.Range("H120").Value = 7494.75
$ws.Range("I120").Value = 7494.75
$ws.Range("K120").Value = 22484.25
$ws.Range("M120").Value = -17646.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 68717.5
$ws.Range("J32").Value = 68717.5
$ws.Range("L32").Value = 68717.5
$ws.Range("N32").Value = -69309.5
$ws.Range("H102").Value = 2938.2222
$ws.Range("I102").Value = 2740.8333
$ws.Range("J102").Value = 3333
$ws.Range("K102").Value = 2740.8333
$ws.Range("L102").Value = 3333
$ws.Range("M102").Value = -1118.8333
$ws.Range("N102").Value = -6577
$ws.Range("H122").Value = 4894.6665
$ws.Range("I122").Value = 4894.6665
$ws.Range("K122").Value = 14683.9995
$ws.Range("M122").Value = -12233.9995
$ws.Range("H126").Value = 4540.222
$ws.Range("I126").Value = 3837.4285
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 11512.2855
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -9042.2855
$ws.Range("N126").Value = -25940
$ws.Range("H132").Value = 5214.1665
$ws.Range("I132").Value = 4932.647
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 14797.941
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -12267.941
$ws.Range("N132").Value = -35060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1396.6666
$ws.Range("I7").Value = 1380.7693
$ws.Range("K7").Value = 1380.7693
$ws.Range("M7").Value = -1268.7693
$ws.Range("H45").Value = 49999
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H48").Value = 39166.25
$ws.Range("I48").Value = 35555.332
$ws.Range("K48").Value = 35555.332
$ws.Range("M48").Value = -34894.332
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H126").Value = 1396.6666
$ws.Range("I126").Value = 1380.7693
$ws.Range("K126").Value = 4142.3079
$ws.Range("M126").Value = -1672.3079
$ws.Range("H132").Value = 9450.333000000001
$ws.Range("I132").Value = 7940.067
$ws.Range("J132").Value = 17001.666
$ws.Range("K132").Value = 23820.201
$ws.Range("L132").Value = 51004.99800000001
$ws.Range("M132").Value = -21290.201
$ws.Range("N132").Value = -56064.99800000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4144.6333
$ws.Range("I107").Value = 3437.5652
$ws.Range("J107").Value = 6467.857
$ws.Range("K107").Value = 10312.6956
$ws.Range("L107").Value = 19403.571
$ws.Range("M107").Value = -8392.695599999999
$ws.Range("N107").Value = -23243.571
$ws.Range("H126").Value = 1930.6207
$ws.Range("I126").Value = 1626.88
$ws.Range("J126").Value = 3829
$ws.Range("K126").Value = 4880.64
$ws.Range("L126").Value = 11487
$ws.Range("M126").Value = -2410.64
$ws.Range("N126").Value = -16427
$ws.Range("H132").Value = 2098.029
$ws.Range("I132").Value = 2076.5085
$ws.Range("J132").Value = 2225
$ws.Range("K132").Value = 6229.5255
$ws.Range("L132").Value = 6675
$ws.Range("M132").Value = -3699.5255
$ws.Range("N132").Value = -11735
